$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Candidate ID 23102882 -> 231102274)
$ws.Range("A2").Value = "uzDwP395"
$ws.Range("B2").Value = 231102274
$ws.Range("C2").Value = "zbvswag56"
$ws.Range("D2").Value = "Pc&9sA8#"
$ws.Range("F2").Value = "TzGfrOCy"
$ws.Range("G2").Value = "XEqY"

# Row 3 (Candidate ID 23102881 -> 231102273)
$ws.Range("A3").Value = "owgaq229"
$ws.Range("B3").Value = 231102273
$ws.Range("C3").Value = "dsjghdk70"
$ws.Range("D3").Value = "E8yBm3&#"
$ws.Range("F3").Value = "puAlNsVg"
$ws.Range("G3").Value = "hXaA"
